# Adds periodic & upfront related scenarios:
# - "repaymentstrategy" (cell B16 on ProductLoanInput) is changed from the
#   placeholder "Mifos style" to "Penalties, Fees, Interest, Principal order",
#   with the cell re-styled to left/top aligned text.
# - ProductLoanInput becomes the active/selected sheet (with B16 selected),
#   replacing ProductLoanOutput as the active tab.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ProductLoanInput")

# Update the repayment-strategy scenario value.
$ws1.Range("B16").Value() = "Penalties, Fees, Interest, Principal order"

# Give it a dedicated left/top aligned style instead of the shared one.
$ws1.Range("B16").HorizontalAlignment = -4131
$ws1.Range("B16").VerticalAlignment = -4160

# Make ProductLoanInput the active sheet/tab, with B16 as the selection.
$ws1.Activate()
$ws1.Range("B16").Select()
